$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("study")
# "study" is also the workbook's active sheet (tabSelected="1"), so this is
# equivalent to `$wb.ActiveSheet`, but naming it explicitly is more robust.

# The "study" sheet gains two new header columns, "title" and "description",
# inserted right after "contact_institution" and before the existing
# "workflow"/"modality" columns. Inserting whole columns at E:F shifts the
# existing E/F ("workflow"/"modality") headers - along with their column
# widths/formatting - out to G/H automatically.
$ws.Range("E:F").Insert()

$ws.Range("E1").Value = "title"
$ws.Range("F1").Value = "description"

# Match the column widths recorded for the new columns.
$ws.Columns.Item(5).ColumnWidth = 3.75
$ws.Columns.Item(6).ColumnWidth = 10.085
